# Re-add multi line char fix to pfrs after conflict (#1474)
#
# 1. Append "f" to the organizationName merge-field (typo/fix re-added).
# 2. Insert ":convCRLF" before ":ifEM():show(.noData)" in a set of
#    merge-fields so multi-line text is converted correctly.

$d = $word.ActiveDocument
$wdReplaceOne = 1

function Replace-Text($find, $replace) {
    $range = $d.Content
    $ok = $range.Find.Execute($find, $true, $false, $false, $false, $false, `
                               $true, 1, $false, $replace, $wdReplaceOne)
    if (-not $ok) {
        throw "Find/Replace failed for: $find"
    }
}

# 1) organizationName field gets a trailing "f"
Replace-Text "{d.parcels[i].owners[i].organizationName}" `
             "{d.parcels[i].owners[i].organizationName}f"

# 2) Add :convCRLF to the following merge-fields
$fields = @(
    "d.parcelsAgricultureDescription",
    "d.parcelsAgricultureImprovementDescription",
    "d.parcelsNonAgricultureUseDescription",
    "d.purpose",
    "d.soilFillTypeToPlace",
    "d.soilTypeRemoved",
    "d.soilStructureFarmUseReason",
    "d.soilStructureResidentialUseReason",
    "d.soilAgriParcelActivity",
    "d.soilStructureResidentialAccessoryUseReason",
    "d.soilStructureOtherUseReason"
)

foreach ($field in $fields) {
    $find = "{" + $field + ":ifEM():show(.noData)}"
    $replace = "{" + $field + ":convCRLF:ifEM():show(.noData)}"
    Replace-Text $find $replace
}

Write-Host "Done"
